# "update to use powerpoint flow"
#
# The author renamed the SE/LL95/UL95 statistics header labels so each one
# carries its outcome prefix (CIN2+ / CIN3+ / CANCER) instead of the bare
# generic label that used to be reused across the three outcome blocks.
# The "Unweighted ..." block headers were similarly normalised so every
# column in that block starts with "Unweighted " (fixing the missing
# space in "UnweightedNumber of ..." along the way).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerUpdates = @{
    'U1'  = 'CIN2+ SE immediate'
    'V1'  = 'CIN2+ LL95 immediate'
    'W1'  = 'CIN2+ UL95 immediate'
    'Y1'  = 'CIN2+ SE 1-year'
    'Z1'  = 'CIN2+ LL95 1-year'
    'AA1' = 'CIN2+ UL95 1-year'
    'AC1' = 'CIN2+ SE 2-year'
    'AD1' = 'CIN2+ LL95 2-year'
    'AE1' = 'CIN2+ UL95 2-year'
    'AG1' = 'CIN2+ SE 3-year'
    'AH1' = 'CIN2+ LL95 3-year'
    'AI1' = 'CIN2+ UL95 3-year'
    'AK1' = 'CIN2+ SE 4-year'
    'AL1' = 'CIN2+ LL95 4-year'
    'AM1' = 'CIN2+ UL95 4-year'
    'AO1' = 'CIN2+ SE 5-year'
    'AP1' = 'CIN2+ LL95 5-year'
    'AQ1' = 'CIN2+ UL95 5-year'
    'AS1' = 'CIN3+ SE immediate'
    'AT1' = 'CIN3+ LL95 immediate'
    'AU1' = 'CIN3+ UL95 immediate'
    'AW1' = 'CIN3+ SE 1-year'
    'AX1' = 'CIN3+ LL95 1-year'
    'AY1' = 'CIN3+ UL95 1-year'
    'BA1' = 'CIN3+ SE 2-year'
    'BB1' = 'CIN3+ LL95 2-year'
    'BC1' = 'CIN3+ UL95 2-year'
    'BE1' = 'CIN3+ SE 3-year'
    'BF1' = 'CIN3+ LL95 3-year'
    'BG1' = 'CIN3+ UL95 3-year'
    'BI1' = 'CIN3+ SE 4-year'
    'BJ1' = 'CIN3+ LL95 4-year'
    'BK1' = 'CIN3+ UL95 4-year'
    'BM1' = 'CIN3+ SE 5-year'
    'BN1' = 'CIN3+ LL95 5-year'
    'BO1' = 'CIN3+ UL95 5-year'
    'BQ1' = 'CANCER SE immediate'
    'BR1' = 'CANCER LL95 immediate'
    'BS1' = 'CANCER UL95 immediate'
    'BU1' = 'CANCER SE 1-year'
    'BV1' = 'CANCER LL95 1-year'
    'BW1' = 'CANCER UL95 1-year'
    'BY1' = 'CANCER SE 2-year'
    'BZ1' = 'CANCER LL95 2-year'
    'CA1' = 'CANCER UL95 2-year'
    'CC1' = 'CANCER SE 3-year'
    'CD1' = 'CANCER LL95 3-year'
    'CE1' = 'CANCER UL95 3-year'
    'CG1' = 'CANCER SE 4-year'
    'CH1' = 'CANCER LL95 4-year'
    'CI1' = 'CANCER UL95 4-year'
    'CK1' = 'CANCER SE 5-year'
    'CL1' = 'CANCER LL95 5-year'
    'CM1' = 'CANCER UL95 5-year'
    'CR1' = 'Unweighted %'
    'CT1' = 'Unweighted CIN2+ Prevalence Cases'
    'CU1' = 'Unweighted CIN2+ Incidence Cases'
    'CV1' = 'Unweighted CIN2+ Unknown Cases'
    'CW1' = 'Unweighted Number of CIN3+ Cases'
    'CX1' = 'Unweighted CIN3+ Prevalence Cases'
    'CY1' = 'Unweighted CIN3+ Incidence Cases'
    'CZ1' = 'Unweighted CIN3+ Unknown Cases'
    'DA1' = 'Unweighted Number of Cancer Cases'
    'DB1' = 'Unweighted Cancer Prevalence Cases'
    'DC1' = 'Unweighted Cancer Incidence Cases'
    'DD1' = 'Unweighted Cancer Unknown Cases'
}

foreach ($addr in $headerUpdates.Keys) {
    $ws.Range($addr).Value2 = $headerUpdates[$addr]
}

# Restore the frozen header pane (columns A:D / row 1 stay pinned) and move
# the visible viewport over to the far right of the newly widened table,
# then leave the whole header row selected as the active selection.
$ws.Range("CY2").Select()
$ws.Rows(1).Select()
